$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cells I1 and J1, copying the style/format from H1
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Fill in I2:J81 with the I0 / IF values (no special style, matching existing data cells)
$data = New-Object 'object[,]' 80,2
$data[0,0] = 5
$data[0,1] = 5
$data[1,0] = 6
$data[1,1] = 6
$data[2,0] = 9
$data[2,1] = 9
$data[3,0] = 8
$data[3,1] = 8
$data[4,0] = 7
$data[4,1] = 7
$data[5,0] = 4
$data[5,1] = 5
$data[6,0] = 5
$data[6,1] = 6
$data[7,0] = 6
$data[7,1] = 6
$data[8,0] = 6
$data[8,1] = 7
$data[9,0] = 6
$data[9,1] = 6
$data[10,0] = 6
$data[10,1] = 6
$data[11,0] = 4
$data[11,1] = 5
$data[12,0] = 4
$data[12,1] = 5
$data[13,0] = 7
$data[13,1] = 7
$data[14,0] = 7
$data[14,1] = 8
$data[15,0] = 6
$data[15,1] = 6
$data[16,0] = 6
$data[16,1] = 6
$data[17,0] = 3
$data[17,1] = 4
$data[18,0] = 7
$data[18,1] = 7
$data[19,0] = 7
$data[19,1] = 7
$data[20,0] = 9
$data[20,1] = 9
$data[21,0] = 5
$data[21,1] = 5
$data[22,0] = 7
$data[22,1] = 7
$data[23,0] = 7
$data[23,1] = 7
$data[24,0] = 5
$data[24,1] = 6
$data[25,0] = 5
$data[25,1] = 6
$data[26,0] = 7
$data[26,1] = 7
$data[27,0] = 10
$data[27,1] = 10
$data[28,0] = 6
$data[28,1] = 7
$data[29,0] = 8
$data[29,1] = 8
$data[30,0] = 9
$data[30,1] = 9
$data[31,0] = 6
$data[31,1] = 6
$data[32,0] = 7
$data[32,1] = 7
$data[33,0] = 8
$data[33,1] = 8
$data[34,0] = 6
$data[34,1] = 6
$data[35,0] = 7
$data[35,1] = 7
$data[36,0] = 9
$data[36,1] = 9
$data[37,0] = 7
$data[37,1] = 7
$data[38,0] = 7
$data[38,1] = 7
$data[39,0] = 5
$data[39,1] = 6
$data[40,0] = 9
$data[40,1] = 9
$data[41,0] = 8
$data[41,1] = 8
$data[42,0] = 8
$data[42,1] = 8
$data[43,0] = 9
$data[43,1] = 9
$data[44,0] = 9
$data[44,1] = 9
$data[45,0] = 8
$data[45,1] = 8
$data[46,0] = 6
$data[46,1] = 7
$data[47,0] = 9
$data[47,1] = 9
$data[48,0] = 7
$data[48,1] = 7
$data[49,0] = 8
$data[49,1] = 8
$data[50,0] = 9
$data[50,1] = 9
$data[51,0] = 8
$data[51,1] = 9
$data[52,0] = 5
$data[52,1] = 6
$data[53,0] = 7
$data[53,1] = 8
$data[54,0] = 6
$data[54,1] = 7
$data[55,0] = 6
$data[55,1] = 6
$data[56,0] = 6
$data[56,1] = 6
$data[57,0] = 9
$data[57,1] = 9
$data[58,0] = 10
$data[58,1] = 10
$data[59,0] = 8
$data[59,1] = 8
$data[60,0] = 3
$data[60,1] = 3
$data[61,0] = 7
$data[61,1] = 7
$data[62,0] = 6
$data[62,1] = 7
$data[63,0] = 7
$data[63,1] = 7
$data[64,0] = 6
$data[64,1] = 6
$data[65,0] = 7
$data[65,1] = 8
$data[66,0] = 11
$data[66,1] = 11
$data[67,0] = 5
$data[67,1] = 5
$data[68,0] = 8
$data[68,1] = 8
$data[69,0] = 4
$data[69,1] = 5
$data[70,0] = 6
$data[70,1] = 6
$data[71,0] = 6
$data[71,1] = 6
$data[72,0] = 7
$data[72,1] = 7
$data[73,0] = 9
$data[73,1] = 9
$data[74,0] = 6
$data[74,1] = 6
$data[75,0] = 6
$data[75,1] = 6
$data[76,0] = 8
$data[76,1] = 8
$data[77,0] = 3
$data[77,1] = 3
$data[78,0] = 3
$data[78,1] = 3
$data[79,0] = 4
$data[79,1] = 4
$ws.Range("I2:J81").Value = $data

Write-Output "Applied I0/IF columns successfully"
